$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.333.53'
$ws.Cells.Item(2, 5).Value = '  +1.49%  '

$ws.Cells.Item(3, 4).Value = '1.905.53'
$ws.Cells.Item(3, 5).Value = '  +1.44%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.002'
$ws.Cells.Item(4, 4).Style = $ws.Cells.Item(4, 3).Style
$ws.Cells.Item(4, 5).Value = '  -0.10%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '323.30'
$ws.Cells.Item(5, 4).Style = $ws.Cells.Item(5, 3).Style
$ws.Cells.Item(5, 5).Value = '  -2.13%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.002'
$ws.Cells.Item(6, 4).Style = $ws.Cells.Item(6, 3).Style
$ws.Cells.Item(6, 5).Value = '  +0.01%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4726'
$ws.Cells.Item(7, 4).Style = $ws.Cells.Item(7, 3).Style

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.4037'
$ws.Cells.Item(8, 4).Style = $ws.Cells.Item(8, 3).Style
$ws.Cells.Item(8, 5).Value = '  -0.60%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.08021'
$ws.Cells.Item(9, 4).Style = $ws.Cells.Item(9, 3).Style
$ws.Cells.Item(9, 5).Value = '  +0.77%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.9937'
$ws.Cells.Item(10, 4).Style = $ws.Cells.Item(10, 3).Style
$ws.Cells.Item(10, 5).Value = '  +0.69%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '22.61'
$ws.Cells.Item(11, 4).Style = $ws.Cells.Item(11, 3).Style
$ws.Cells.Item(11, 5).Value = '  +4.90%  '

$ws.Cells.Item(12, 4).Value = '1.905.63'
$ws.Cells.Item(12, 5).Value = '  +1.45%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '5.857'
$ws.Cells.Item(13, 4).Style = $ws.Cells.Item(13, 3).Style
$ws.Cells.Item(13, 5).Value = '  -0.52%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '7.054'
$ws.Cells.Item(14, 4).Style = $ws.Cells.Item(14, 3).Style
$ws.Cells.Item(14, 5).Value = '  +0.04%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '89.32'
$ws.Cells.Item(15, 4).Style = $ws.Cells.Item(15, 3).Style
$ws.Cells.Item(15, 5).Value = '  +1.42%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '1.002'
$ws.Cells.Item(16, 4).Style = $ws.Cells.Item(16, 3).Style
$ws.Cells.Item(16, 5).Value = '  -0.15%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.06630'
$ws.Cells.Item(17, 4).Style = $ws.Cells.Item(17, 3).Style
$ws.Cells.Item(17, 5).Value = '  +1.08%  '

$ws.Cells.Item(18, 5).Value = '  +0.06%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '17.55'
$ws.Cells.Item(19, 4).Style = $ws.Cells.Item(19, 3).Style
$ws.Cells.Item(19, 5).Value = '  +1.05%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '1.002'
$ws.Cells.Item(20, 4).Style = $ws.Cells.Item(20, 3).Style
$ws.Cells.Item(20, 5).Value = '  -0.03%  '

$ws.Cells.Item(21, 4).Value = '29.336.37'
$ws.Cells.Item(21, 5).Value = '  +1.41%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.510'
$ws.Cells.Item(22, 4).Style = $ws.Cells.Item(22, 3).Style
$ws.Cells.Item(22, 5).Value = '  +2.04%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '11.44'
$ws.Cells.Item(23, 4).Style = $ws.Cells.Item(23, 3).Style
$ws.Cells.Item(23, 5).Value = '  -0.25%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.204'
$ws.Cells.Item(24, 4).Style = $ws.Cells.Item(24, 3).Style
$ws.Cells.Item(24, 5).Value = '  +0.13%  '

$ws.Cells.Item(25, 4).Value = '2.117.12'
$ws.Cells.Item(25, 5).Value = '  +0.36%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '154.18'
$ws.Cells.Item(26, 4).Style = $ws.Cells.Item(26, 3).Style
$ws.Cells.Item(26, 5).Value = '  -1.46%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '19.74'
$ws.Cells.Item(27, 4).Style = $ws.Cells.Item(27, 3).Style
$ws.Cells.Item(27, 5).Value = '  +1.40%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '6.053'
$ws.Cells.Item(28, 4).Style = $ws.Cells.Item(28, 3).Style
$ws.Cells.Item(28, 5).Value = '  +11.07%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.089'
$ws.Cells.Item(29, 4).Style = $ws.Cells.Item(29, 3).Style
$ws.Cells.Item(29, 5).Value = '  +0.81%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '117.60'
$ws.Cells.Item(30, 4).Style = $ws.Cells.Item(30, 3).Style
$ws.Cells.Item(30, 5).Value = '  +0.33%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.068'
$ws.Cells.Item(31, 4).Style = $ws.Cells.Item(31, 3).Style
$ws.Cells.Item(31, 5).Value = '  +4.47%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.09490'
$ws.Cells.Item(32, 4).Style = $ws.Cells.Item(32, 3).Style
$ws.Cells.Item(32, 5).Value = '  +1.78%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.413'
$ws.Cells.Item(33, 4).Style = $ws.Cells.Item(33, 3).Style
$ws.Cells.Item(33, 5).Value = '  +1.27%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '3.558'
$ws.Cells.Item(34, 4).Style = $ws.Cells.Item(34, 3).Style
$ws.Cells.Item(34, 5).Value = '  +2.16%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '5.356'
$ws.Cells.Item(35, 4).Style = $ws.Cells.Item(35, 3).Style
$ws.Cells.Item(35, 5).Value = '  +1.76%  '

$ws.Cells.Item(36, 2).Value = 'Hedera'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.06068'
$ws.Cells.Item(36, 4).Style = $ws.Cells.Item(36, 3).Style
$ws.Cells.Item(36, 5).Value = '  +0.71%  '

$ws.Cells.Item(37, 2).Value = 'VeChain'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.02247'
$ws.Cells.Item(37, 4).Style = $ws.Cells.Item(37, 3).Style
$ws.Cells.Item(37, 5).Value = '  +1.23%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '1.173'
$ws.Cells.Item(38, 4).Style = $ws.Cells.Item(38, 3).Style
$ws.Cells.Item(38, 5).Value = '  +0.18%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '8.044'
$ws.Cells.Item(39, 4).Style = $ws.Cells.Item(39, 3).Style
$ws.Cells.Item(39, 5).Value = '  -2.74%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.5811'
$ws.Cells.Item(40, 4).Style = $ws.Cells.Item(40, 3).Style
$ws.Cells.Item(40, 5).Value = '  +0.82%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '2.501'
$ws.Cells.Item(41, 4).Style = $ws.Cells.Item(41, 3).Style
$ws.Cells.Item(41, 5).Value = '  +11.09%  '

$ws.Cells.Item(42, 5).Value = '  +0.65%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '10.07'
$ws.Cells.Item(43, 4).Style = $ws.Cells.Item(43, 3).Style
$ws.Cells.Item(43, 5).Value = '  +0.24%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.07800'
$ws.Cells.Item(44, 4).Style = $ws.Cells.Item(44, 3).Style
$ws.Cells.Item(44, 5).Value = '  +4.25%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.273'
$ws.Cells.Item(45, 4).Style = $ws.Cells.Item(45, 3).Style
$ws.Cells.Item(45, 5).Value = '  +1.43%  '

$ws.Cells.Item(46, 5).Value = '  +1.10%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.5483'
$ws.Cells.Item(47, 4).Style = $ws.Cells.Item(47, 3).Style
$ws.Cells.Item(47, 5).Value = '  +0.90%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.903'
$ws.Cells.Item(48, 4).Style = $ws.Cells.Item(48, 3).Style
$ws.Cells.Item(48, 5).Value = '  +0.44%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '113.40'
$ws.Cells.Item(49, 4).Style = $ws.Cells.Item(49, 3).Style
$ws.Cells.Item(49, 5).Value = '  +2.22%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '43.80'
$ws.Cells.Item(50, 4).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(50, 5).Value = '  -3.14%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.2920'
$ws.Cells.Item(51, 4).Style = $ws.Cells.Item(51, 3).Style
$ws.Cells.Item(51, 5).Value = '  +3.72%  '
